$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Remove the "Side note on general:" block (heading + three git-command
# paragraphs) that had been appended after "Resolve any merge conflicts
# if they arise.". Locate paragraphs by their text (robust to index
# drift) rather than hard-coded paragraph numbers.
# ---------------------------------------------------------------------

$count = $d.Paragraphs.Count
$resolveIdx = -1
$lastPushIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Resolve any merge conflicts*") {
        $resolveIdx = $i
    }
    if ($t -like "*git push*") {
        $lastPushIdx = $i
    }
}

if ($resolveIdx -gt 0 -and $lastPushIdx -gt $resolveIdx) {
    # Delete from the start of the paragraph right after "Resolve any merge
    # conflicts if they arise." through the end of the trailing "git push"
    # paragraph. That span covers the blank paragraph, the "Side note on
    # general:" heading, the "Run these three commands" paragraph, the git
    # commit-message paragraph, and the final git push paragraph -- leaving
    # the document ending in a single empty paragraph, same as before the
    # note was appended.
    $delStart = $d.Paragraphs.Item($resolveIdx + 1).Range.Start
    $delEnd = $d.Paragraphs.Item($lastPushIdx).Range.End
    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}

# Safety net: if two empty paragraphs ended up adjacent right after
# "Resolve any merge conflicts if they arise.", collapse them into one.
# (Paragraph.Range.Text includes the trailing paragraph mark, so trim it
# before checking for "no visible text".)
$afterIdx = $resolveIdx + 1
if ($afterIdx -lt $d.Paragraphs.Count) {
    $cur = $d.Paragraphs.Item($afterIdx)
    $nxt = $d.Paragraphs.Item($afterIdx + 1)
    if ($cur.Range.Text.Trim() -eq "" -and $nxt.Range.Text.Trim() -eq "") {
        $cur.Range.Delete()
    }
}

# ---------------------------------------------------------------------
# Remove the now-unused "HTML Preformatted" paragraph/character styles.
# Grab both style references before deleting either (deleting the base
# style first while still holding the linked char-style reference can
# otherwise leave the char style behind).
# ---------------------------------------------------------------------
$preChar = $d.Styles.Item("HTML Preformatted Char")
$pre = $d.Styles.Item("HTML Preformatted")
$preChar.Delete()
$pre.Delete()
